$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename header row on "Intake" sheet (columns A-G), keeping column order,
# just collapsing the old spaced/hyphenated labels into compact identifiers.
$ws1.Range("A1").Value = "AcademicYear"
$ws1.Range("B1").Value = "AssociateDegree"
$ws1.Range("C1").Value = "HigherDiploma"
$ws1.Range("D1").Value = "Subdegree"
$ws1.Range("E1").Value = "FirstYearFirstDegree"
$ws1.Range("G1").Value = "Undergraduate"
$ws1.Range("F1").Value = "TopUpDegree"

# Same header rename on "Graduates" sheet.
$ws2.Range("A1").Value = "AcademicYear"
$ws2.Range("B1").Value = "AssociateDegree"
$ws2.Range("C1").Value = "HigherDiploma"
$ws2.Range("D1").Value = "Subdegree"
$ws2.Range("E1").Value = "FirstYearFirstDegree"
$ws2.Range("G1").Value = "Undergraduate"
$ws2.Range("F1").Value = "TopUpDegree"

# Auto-fit the Graduates columns to their (now different) header text.
$ws2.Columns.Item(1).ColumnWidth = 12.333333333333332
$ws2.Columns.Item(2).ColumnWidth = 14.5
$ws2.Columns.Item(3).ColumnWidth = 13.0
$ws2.Columns.Item(5).ColumnWidth = 18.166666666666668
$ws2.Columns.Item(6).ColumnWidth = 12.166666666666666
$ws2.Columns.Item(7).ColumnWidth = 13.333333333333332

# Select the header row on both sheets, and make "Intake" the active tab.
$ws2.Range("A1:G1").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A1:G1").Select() | Out-Null
